$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.599.11'
$ws.Range("E2").Value = '  +6.11%  '
$ws.Range("D3").Value = '3.577.64'
$ws.Range("E3").Value = '  +5.47%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '''591.22'
$ws.Range("E5").Value = '  +5.64%  '
$ws.Range("D6").Value = '''190.88'
$ws.Range("E6").Value = '  +8.83%  '
$ws.Range("D7").Value = '''0.647'
$ws.Range("E7").Value = '  +2.68%  '
$ws.Range("D8").Value = '3.573.64'
$ws.Range("E8").Value = '  +5.65%  '
$ws.Range("E9").Value = '  -0.10%  '
$ws.Range("D10").Value = '''0.181'
$ws.Range("E10").Value = '  +4.07%  '
$ws.Range("D11").Value = '''0.660'
$ws.Range("E11").Value = '  +4.15%  '
$ws.Range("D12").Value = '''57.86'
$ws.Range("E12").Value = '  +7.98%  '
$ws.Range("D13").Value = '''0.0000293'
$ws.Range("E13").Value = '  +5.41%  '
$ws.Range("D14").Value = '''9.70'
$ws.Range("E14").Value = '  +5.54%  '
$ws.Range("D15").Value = '4.155.85'
$ws.Range("E15").Value = '  +5.60%  '
$ws.Range("D16").Value = '''19.32'
$ws.Range("E16").Value = '  +5.90%  '
$ws.Range("D17").Value = '3.577.98'
$ws.Range("E17").Value = '  +5.46%  '
$ws.Range("D18").Value = '69.609.64'
$ws.Range("E18").Value = '  +6.01%  '
$ws.Range("D19").Value = '''12.64'
$ws.Range("E19").Value = '  +6.77%  '
$ws.Range("E20").Value = '  +1.11%  '
$ws.Range("E21").Value = '  +4.95%  '
$ws.Range("D22").Value = '''499.81'
$ws.Range("E22").Value = '  +3.75%  '
$ws.Range("E23").Value = '  +10.13%  '
$ws.Range("D24").Value = '''16.89'
$ws.Range("E24").Value = '  +18.15%  '
$ws.Range("D25").Value = '''4.45'
$ws.Range("E25").Value = '  +8.80%  '
$ws.Range("D26").Value = '''90.94'
$ws.Range("E26").Value = '  +0.99%  '
$ws.Range("D27").Value = '''3.08'
$ws.Range("E27").Value = '  +5.67%  '
$ws.Range("D28").Value = '''11.10'
$ws.Range("E28").Value = '  +4.67%  '
$ws.Range("D29").Value = '''9.33'
$ws.Range("E29").Value = '  +6.96%  '
$ws.Range("D30").Value = '''32.17'
$ws.Range("E30").Value = '  +2.75%  '
$ws.Range("D31").Value = '''7.49'
$ws.Range("E31").Value = '  +14.09%  '
$ws.Range("D32").Value = '''12.15'
$ws.Range("E32").Value = '  +6.27%  '
$ws.Range("D33").Value = '''614.24'
$ws.Range("E33").Value = '  +7.33%  '
$ws.Range("D34").Value = '''65.39'
$ws.Range("E34").Value = '  +2.71%  '
$ws.Range("E35").Value = '  +7.36%  '
$ws.Range("D36").Value = '0.0₃0828'
$ws.Range("E36").Value = '  +11.33%  '
$ws.Range("B37").Value = 'Dai'
$ws.Range("C37").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D37").Value = '''1.00'
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").Value = '''0.147'
$ws.Range("E38").Value = '  +3.91%  '
$ws.Range("D39").Value = '''37.86'
$ws.Range("E39").Value = '  +5.78%  '
$ws.Range("E40").Value = '  +6.60%  '
$ws.Range("D41").Value = '''3.62'
$ws.Range("E41").Value = '  -0.68%  '
$ws.Range("D42").Value = '3.340.72'
$ws.Range("E42").Value = '  +8.10%  '
$ws.Range("E43").Value = '  +10.96%  '
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = '''0.0443'
$ws.Range("E44").Value = '  +6.63%  '
$ws.Range("B45").Value = 'Fetch.AI'
$ws.Range("C45").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D45").Value = '''2.68'
$ws.Range("E45").Value = '  +9.66%  '
$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").Value = '''0.138'
$ws.Range("E46").Value = '  +2.86%  '
$ws.Range("D47").Value = '''3.27'
$ws.Range("E47").Value = '  +3.08%  '
$ws.Range("B48").Value = 'dogwifhat'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D48").Value = '''2.78'
$ws.Range("E48").Value = '  +14.96%  '
$ws.Range("D49").Value = '''9.05'
$ws.Range("E49").Value = '  +7.40%  '
$ws.Range("B50").Value = 'FirstDigitalUSD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D50").Value = '''1.00'
$ws.Range("E50").Value = '  +0.25%  '
$ws.Range("B51").Value = 'LidoDAOToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D51").Value = '''3.23'
$ws.Range("E51").Value = '  +4.67%  '
